$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the values currently held by the merged (top-left) cells ---
$valA2  = $ws.Range("A2").Value2
$valA18 = $ws.Range("A18").Value2
$valB2  = $ws.Range("B2").Value2
$valB8  = $ws.Range("B8").Value2
$valB12 = $ws.Range("B12").Value2
$valB18 = $ws.Range("B18").Value2
$valB23 = $ws.Range("B23").Value2
$valB26 = $ws.Range("B26").Value2

# --- Unmerge every merged block in columns A and B ---
$ws.Range("A2:A17").UnMerge()
$ws.Range("A18:A32").UnMerge()
$ws.Range("B2:B7").UnMerge()
$ws.Range("B8:B11").UnMerge()
$ws.Range("B12:B17").UnMerge()
$ws.Range("B18:B22").UnMerge()
$ws.Range("B23:B25").UnMerge()
$ws.Range("B26:B32").UnMerge()

# --- Fill every row of each former merge block with its value ---
# Column A: first block (OVERHEAD) gets the value repeated on every row;
# the second block (LAVATORY) keeps the value only on its first row.
$ws.Range("A2:A17").Value = $valA2
$ws.Range("A18").Value = $valA18

$ws.Range("B2:B7").Value = $valB2
$ws.Range("B8:B11").Value = $valB8
$ws.Range("B12:B17").Value = $valB12
$ws.Range("B18:B22").Value = $valB18
$ws.Range("B23:B25").Value = $valB23
$ws.Range("B26:B32").Value = $valB26

# --- Re-apply the plain bordered / centred format (same as column C) to A2:A32 and B2:B32 ---
$ws.Range("C2").Copy()
$ws.Range("A2:A32").PasteSpecial(-4122)
$ws.Range("B2:B32").PasteSpecial(-4122)

# --- Update the active selection shown in the saved view ---
$ws.Range("B27:B32").Select()
